$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B-column (missing fraction) values for rows 2-78
$ws.Range("B2").Value = 0.222582738481506
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0.636534717715769
$ws.Range("B5").Value = 0.636534717715769
$ws.Range("B6").Value = 0.636534717715769
$ws.Range("B7").Value = 0.636534717715769
$ws.Range("B8").Value = 0.636534717715769
$ws.Range("B9").Value = 0.636534717715769
$ws.Range("B10").Value = 0.636534717715769
$ws.Range("B11").Value = 0.636534717715769
$ws.Range("B12").Value = 0.636534717715769
$ws.Range("B13").Value = 0.636534717715769
$ws.Range("B14").Value = 0.636534717715769
$ws.Range("B15").Value = 0.636534717715769
$ws.Range("B16").Value = 0.636534717715769
$ws.Range("B17").Value = 0.636534717715769
$ws.Range("B18").Value = 0.636534717715769
$ws.Range("B19").Value = 0.636534717715769
$ws.Range("B20").Value = 0.636534717715769
$ws.Range("B21").Value = 0.636534717715769
$ws.Range("B22").Value = 0.636534717715769
$ws.Range("B23").Value = 0.636534717715769
$ws.Range("B24").Value = 0.636534717715769
$ws.Range("B25").Value = 0.636534717715769
$ws.Range("B26").Value = 0.636534717715769
$ws.Range("B27").Value = 0.636534717715769
$ws.Range("B28").Value = 0.636534717715769
$ws.Range("B29").Value = 0.636534717715769
$ws.Range("B30").Value = 0.636534717715769
$ws.Range("B31").Value = 0.762773091066407
$ws.Range("B32").Value = 0.763140817650876
$ws.Range("B33").Value = 0.78537746052347
$ws.Range("B34").Value = 0.806727233398226
$ws.Range("B35").Value = 0.774042829331603
$ws.Range("B36").Value = 0.869370538611291
$ws.Range("B37").Value = 0.965801427644387
$ws.Range("B38").Value = 0.965801427644387
$ws.Range("B39").Value = 0.965801427644387
$ws.Range("B40").Value = 0.77255029201817
$ws.Range("B41").Value = 0.772809863724854
$ws.Range("B42").Value = 0.772831494700411
$ws.Range("B43").Value = 0.772809863724854
$ws.Range("B44").Value = 0.829699329439758
$ws.Range("B45").Value = 0.772809863724854
$ws.Range("B46").Value = 0.974129353233831
$ws.Range("B47").Value = 0.994570625135194
$ws.Range("B48").Value = 0.983928185161151
$ws.Range("B49").Value = 0.983928185161151
$ws.Range("B50").Value = 0.983928185161151
$ws.Range("B51").Value = 0.994029850746269
$ws.Range("B52").Value = 0.993813540990699
$ws.Range("B53").Value = 0.994116374648497
$ws.Range("B54").Value = 0.995695435864157
$ws.Range("B55").Value = 0.995695435864157
$ws.Range("B56").Value = 0.995695435864157
$ws.Range("B57").Value = 0.992991563919533
$ws.Range("B58").Value = 0.993207873675103
$ws.Range("B59").Value = 0.993640493186243
$ws.Range("B60").Value = 0.992018170019468
$ws.Range("B61").Value = 0.992710361237292
$ws.Range("B62").Value = 0.993099718797318
$ws.Range("B63").Value = 0.993489076357344
$ws.Range("B64").Value = 0.992883409041748
$ws.Range("B65").Value = 0.99301319489509
$ws.Range("B66").Value = 0.98695652173913
$ws.Range("B67").Value = 0.995695435864157
$ws.Range("B68").Value = 0.995695435864157
$ws.Range("B69").Value = 0.993121349772875
$ws.Range("B70").Value = 0.993813540990699
$ws.Range("B71").Value = 0.993034825870647
$ws.Range("B72").Value = 0.994116374648497
$ws.Range("B73").Value = 0.992731992212849
$ws.Range("B74").Value = 0.992299372701709
$ws.Range("B75").Value = 0.993207873675103
$ws.Range("B76").Value = 0.993207873675103
$ws.Range("B77").Value = 0.993597231235129
$ws.Range("A78").Value = "heq"
$ws.Range("B78").Value = 0.76707765520225
$ws.Range("A79").Value = "gini"
$ws.Range("B79").Value = 0.952563270603504
$ws.Range("A80").Value = "population"
$ws.Range("B80").Value = 0.762773091066407
$ws.Range("A81").Value = "median_age"
$ws.Range("B81").Value = 0.625567813108371
$ws.Range("A82").Value = "aged65"
$ws.Range("B82").Value = 0.697252866104261
$ws.Range("A83").Value = "aged70_male"
$ws.Range("B83").Value = 0.697252866104261
$ws.Range("A84").Value = "aged70_female"
$ws.Range("B84").Value = 0.697252866104261
$ws.Range("A85").Value = "ext_poverty"
$ws.Range("B85").Value = 0.950703006705602
$ws.Range("A86").Value = "life_expectancy"
$ws.Range("B86").Value = 0.714233181916504
$ws.Range("A87").Value = "schooling"
$ws.Range("B87").Value = 0.863962794722042
$ws.Range("A88").Value = "literacy"
$ws.Range("B88").Value = 0.974432186891629
$ws.Range("A89").Value = "urban"
$ws.Range("B89").Value = 0.698615617564352
$ws.Range("A90").Value = "dependency"
$ws.Range("B90").Value = 0.697901795370971
$ws.Range("A91").Value = "homicide"
$ws.Range("B91").Value = 0.90893359290504
$ws.Range("A92").Value = "net_migration"
$ws.Range("B92").Value = 0.70194678780013
$ws.Range("A93").Value = "hdi"
$ws.Range("B93").Value = 0.871533636166991
